$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "High Priority break-up" sheet handling
#    - Rename existing sheet to "Interannual update - High Pri" and replace
#      its data with the new interannual-update numbers.
#    - Add a brand new sheet "Major update - High Priority " placed right
#      after it, holding the data that used to live in the old sheet
#      (i.e. a straight copy of the original content).
# ---------------------------------------------------------------------------

$wsOld = $wb.Worksheets.Item("High Priority break-up")

# Create the new sheet first (as an exact copy) so it keeps the original
# "High Priority break-up" numbers before we overwrite $wsOld with the new
# interannual data.
$wsOld.Copy([System.Reflection.Missing]::Value, $wsOld)
$wsMajor = $wb.Worksheets.Item($wsOld.Index + 1)
$wsMajor.Name = "Major update - High Priority "

# Now rename the original sheet and overwrite its contents.
$wsOld.Name = "Interannual update - High Pri"

$wsInter = $wsOld

$wsInter.Cells.Clear()

$wsInter.Range("A1").Value = "Break-up"
$wsInter.Range("B1").Value = "High Species (no.)"
$wsInter.Range("C1").Value = "High Species (perc.)"
$wsInter.Range("D1").Value = "New High Species (no.)"
$wsInter.Range("E1").Value = "New High Species (perc.)"
$wsInter.Range("A1:E1").Font.Bold = $true
$wsInter.Range("A1:E1").HorizontalAlignment = -4108

$wsInter.Range("A2").Value = "Trend New"
$wsInter.Range("B2").Value = 81
$wsInter.Range("C2").Value = 78.59999999999999
$wsInter.Range("D2").Value = 81
$wsInter.Range("E2").Value = 83.5

$wsInter.Range("A3").Value = "IUCN"
$wsInter.Range("B3").Value = 22
$wsInter.Range("C3").Value = 21.4
$wsInter.Range("D3").Value = 16
$wsInter.Range("E3").Value = 16.5

# ---------------------------------------------------------------------------
# 2. "Trends Status" sheet: update numbers, drop column D for rows 2-6.
# ---------------------------------------------------------------------------

$wsTrends = $wb.Worksheets.Item("Trends Status")

$wsTrends.Range("C2").Value = 1
$wsTrends.Range("D2").ClearContents()
$wsTrends.Range("E2").Value = 50

$wsTrends.Range("C3").Value = 0
$wsTrends.Range("D3").ClearContents()
$wsTrends.Range("E3").Value = 0

$wsTrends.Range("B4").Value = 0
$wsTrends.Range("C4").Value = 1
$wsTrends.Range("D4").ClearContents()
$wsTrends.Range("E4").Value = 50

$wsTrends.Range("D5").ClearContents()

$wsTrends.Range("D6").ClearContents()

$wsTrends.Range("B7").Value = 7
$wsTrends.Range("C7").Value = 28

$wsTrends.Range("B8").Value = 415
$wsTrends.Range("C8").Value = 392

# ---------------------------------------------------------------------------
# 3. "Priority Status" sheet: update species counts.
# ---------------------------------------------------------------------------

$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 4. "Species qualification" sheet: update label and counts.
# ---------------------------------------------------------------------------

$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("A2").Value = "SoIB Assessment"
$wsQual.Range("B2").Value = 422
$wsQual.Range("B3").Value = 7
$wsQual.Range("C3").Value = 0
$wsQual.Range("C4").Value = 2

# ---------------------------------------------------------------------------
# Restore the originally active sheet/tab (sheet-creation/copy operations
# above shift the active tab), so the workbook-level view state is left
# untouched, matching the source workbook.
# ---------------------------------------------------------------------------
$wsTrends.Activate()

$wb.Save()
